# Delete (Fuh, 2000) data from VEGF:NRP1 data
#
# The "VEGFA165_NRP1" worksheet has a row for the "Fuh et al., 2000"
# reference (row 6: Reference | SPR | 113 | =C6*0.35, with a comment
# "Error: ~35%" on D6). This row is removed entirely, which shifts the
# remaining rows up, updates the shared-string table, and drops the
# now-orphaned cell comment.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEGFA165_NRP1")

# Make this the active/selected sheet (it becomes the active tab after
# the edit).
$ws.Activate()

# Remove the comment attached to D6 (tied to the Fuh et al., 2000 row)
# before the row holding it is deleted.
$comment = $ws.Range("D6").Comment
if ($comment -ne $null) {
    $comment.Delete()
}

# Delete the entire row for the Fuh et al., 2000 reference; this shifts
# rows 7-8 up to become rows 6-7.
$ws.Rows.Item(6).Delete()

# Match the resulting selection left behind in the file.
[void]$ws.Range("A6:D6").Select()
